$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '80.870.69'
$ws.Range('E2').Value = '  +2.22%  '

# Row 3
$ws.Range('D3').Value = '3.131.94'
$ws.Range('E3').Value = '  -1.33%  '

# Row 4
$ws.Range('E4').Value = '  +0.10%  '

# Row 5
$ws.Range('D5').Value = '''205.73'
$ws.Range('E5').Value = '  -0.41%  '

# Row 6
$ws.Range('D6').Value = '''617.48'
$ws.Range('E6').Value = '  -1.84%  '

# Row 7
$ws.Range('D7').Value = '''0.280'
$ws.Range('E7').Value = '  +23.31%  '

# Row 8
$ws.Range('E8').Value = '  -0.04%  '

# Row 9
$ws.Range('D9').Value = '''0.575'
$ws.Range('E9').Value = '  -1.70%  '

# Row 10
$ws.Range('D10').Value = '3.132.37'
$ws.Range('E10').Value = '  -1.37%  '

# Row 11
$ws.Range('D11').Value = '''0.572'
$ws.Range('E11').Value = '  -0.33%  '

# Row 12
$ws.Range('D12').Value = '''0.0000249'
$ws.Range('E12').Value = '  +10.91%  '

# Row 13
$ws.Range('E13').Value = '  +0.27%  '

# Row 14
$ws.Range('E14').Value = '  -3.48%  '

# Row 15
$ws.Range('D15').Value = '3.708.54'
$ws.Range('E15').Value = '  -1.21%  '

# Row 16
$ws.Range('D16').Value = '''30.97'
$ws.Range('E16').Value = '  -2.32%  '

# Row 17
$ws.Range('D17').Value = '80.933.64'
$ws.Range('E17').Value = '  +2.51%  '

# Row 18
$ws.Range('D18').Value = '3.130.19'
$ws.Range('E18').Value = '  -1.55%  '

# Row 19
$ws.Range('D19').Value = '''3.13'
$ws.Range('E19').Value = '  +10.38%  '

# Row 20
$ws.Range('D20').Value = '''13.71'
$ws.Range('E20').Value = '  -5.11%  '

# Row 21
$ws.Range('D21').Value = '''425.59'
$ws.Range('E21').Value = '  -0.52%  '

# Row 22
$ws.Range('D22').Value = '''8.84'
$ws.Range('E22').Value = '  -5.93%  '

# Row 23
$ws.Range('D23').Value = '''5.02'
$ws.Range('E23').Value = '  +1.08%  '

# Row 24
$ws.Range('D24').Value = '''7.14'
$ws.Range('E24').Value = '  +5.46%  '

# Row 25
$ws.Range('D25').Value = '''5.12'
$ws.Range('E25').Value = '  +7.68%  '

# Row 26
$ws.Range('D26').Value = '3.303.27'
$ws.Range('E26').Value = '  -1.16%  '

# Row 27
$ws.Range('D27').Value = '''75.21'
$ws.Range('E27').Value = '  -1.84%  '

# Row 28
$ws.Range('D28').Value = '''10.66'
$ws.Range('E28').Value = '  -2.46%  '

# Row 29
$ws.Range('E29').Value = '  -0.20%  '

# Row 30
$ws.Range('D30').Value = '''0.0000119'
$ws.Range('E30').Value = '  +4.55%  '

# Row 31
$ws.Range('E31').Value = '  +0.06%  '

# Row 32
$ws.Range('D32').Value = '''8.88'
$ws.Range('E32').Value = '  +0.08%  '

# Row 33
$ws.Range('D33').Value = '''553.61'
$ws.Range('E33').Value = '  +8.09%  '

# Row 34
$ws.Range('B34').Value = 'Cronos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D34').Value = '''0.149'
$ws.Range('E34').Value = '  +13.89%  '

# Row 35
$ws.Range('B35').Value = 'Fetch.AI'
$ws.Range('C35').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D35').Value = '''1.46'
$ws.Range('E35').Value = '  -0.61%  '

# Row 36
$ws.Range('D36').Value = '''0.149'
$ws.Range('E36').Value = '  +11.02%  '

# Row 37
$ws.Range('D37').Value = '''1.96'
$ws.Range('E37').Value = '  -0.12%  '

# Row 38
$ws.Range('D38').Value = '''22.43'
$ws.Range('E38').Value = '  -1.68%  '

# Row 39
$ws.Range('E39').Value = '  -0.05%  '

# Row 40
$ws.Range('D40').Value = '''0.401'
$ws.Range('E40').Value = '  +0.60%  '

# Row 41
$ws.Range('D41').Value = '''20.70'
$ws.Range('E41').Value = '  +3.65%  '

# Row 42
$ws.Range('D42').Value = '''5.84'
$ws.Range('E42').Value = '  +8.31%  '

# Row 43
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').Value = '''1.97'
$ws.Range('E43').Value = '  +11.68%  '

# Row 44
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').Value = '''2.95'
$ws.Range('E44').Value = '  +17.62%  '

# Row 45
$ws.Range('B45').Value = 'Monero'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D45').Value = '''158.94'
$ws.Range('E45').Value = '  -2.87%  '

# Row 46
$ws.Range('E46').Value = '  +0.04%  '

# Row 47
$ws.Range('D47').Value = '''185.12'
$ws.Range('E47').Value = '  -3.96%  '

# Row 48
$ws.Range('D48').Value = '''44.11'
$ws.Range('E48').Value = '  +3.58%  '

# Row 49
$ws.Range('E49').Value = '  +0.79%  '

# Row 50
$ws.Range('D50').Value = '''0.762'
$ws.Range('E50').Value = '  -6.34%  '

# Row 51
$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').Value = '''25.27'
$ws.Range('E51').Value = '  +3.33%  '
